# Applies the "Some updates to texts" commit:
#  - Bumps the cached datetimeFigureOut placeholder text on the Slide
#    Master and every Slide Layout from 11/16/2020 to 3/18/2021.
#  - Bumps the literal date text run on every slide's Date placeholder
#    from 2020-11-16 to 2021-03-18.
#  - Bumps the copyright year in every slide's Footer placeholder from
#    "2012-2020" to "2012-2021".
#  - Merges the three runs of "Runs AJAX request " / "and receives " /
#    "responses" on slide 3 into a single run of text.

$p = $ppt.ActivePresentation

function Update-PlaceholderText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            $t = $tr.Text
            if ($t -eq "11/16/2020") {
                $tr.Text = "3/18/2021"
            }
        }
    }
}

# Slide Master date placeholder.
Update-PlaceholderText $p.SlideMaster.Shapes

# Every Slide Layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-PlaceholderText $layouts.Item($li).Shapes
}

# Per-slide footer date / copyright text runs.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            $t = $tr.Text
            if ($t -eq "2020-11-16") {
                $tr.Text = "2021-03-18"
            } elseif ($t -match "^.* Juhani .*2012-2020$") {
                $tr.Text = $t -replace "2012-2020$", "2012-2021"
            }
        }
    }
}

# Slide 3: merge "Runs AJAX request " + "and receives " + "responses"
# into a single run's text.
$slide3 = $p.Slides.Item(3)
$shape1 = $slide3.Shapes.Item(1)
$tr = $shape1.TextFrame.TextRange
$full = $tr.Text
$marker = "Runs AJAX request "
$idx = $full.IndexOf($marker)
if ($idx -ge 0) {
    $oldPhrase = "Runs AJAX request and receives responses"
    $sub = $tr.Characters($idx + 1, $oldPhrase.Length)
    $sub.Text = $oldPhrase
}
